$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell address ("D2", "E2", ...) from the crypto price/volume
# table and the new literal text that should be displayed there. The sheet
# stores every Price/Volume value as text (inline strings), so values that
# look like plain numbers ("301.94") are written with a leading apostrophe
# to keep Excel from re-interpreting them as numeric, then the cell style is
# reset back to Normal so no visible/Text-format styling is left behind.
$updates = @(
    @{ Cell = "D2"; Value = '42.728.68' },
    @{ Cell = "E2"; Value = '  -1.15%  ' },
    @{ Cell = "D3"; Value = '2.317.72' },
    @{ Cell = "E3"; Value = '  -0.31%  ' },
    @{ Cell = "E4"; Value = '  +0.10%  ' },
    @{ Cell = "D5"; Value = '301.94' },
    @{ Cell = "E5"; Value = '  -0.31%  ' },
    @{ Cell = "D6"; Value = '93.54' },
    @{ Cell = "E6"; Value = '  -4.26%  ' },
    @{ Cell = "D7"; Value = '0.499' },
    @{ Cell = "E7"; Value = '  -1.31%  ' },
    @{ Cell = "E8"; Value = '  +0.14%  ' },
    @{ Cell = "D9"; Value = '0.490' },
    @{ Cell = "E9"; Value = '  -2.32%  ' },
    @{ Cell = "D10"; Value = '33.84' },
    @{ Cell = "E10"; Value = '  -4.99%  ' },
    @{ Cell = "D11"; Value = '0.0778' },
    @{ Cell = "E11"; Value = '  -2.46%  ' },
    @{ Cell = "D12"; Value = '18.52' },
    @{ Cell = "E12"; Value = '  -5.49%  ' },
    @{ Cell = "D13"; Value = '0.121' },
    @{ Cell = "E13"; Value = '  +1.19%  ' },
    @{ Cell = "D14"; Value = '6.66' },
    @{ Cell = "E14"; Value = '  -4.01%  ' },
    @{ Cell = "D15"; Value = '2.690.81' },
    @{ Cell = "E15"; Value = '  +0.22%  ' },
    @{ Cell = "D16"; Value = '2.336.34' },
    @{ Cell = "E16"; Value = '  +0.35%  ' },
    @{ Cell = "D17"; Value = '0.786' },
    @{ Cell = "E17"; Value = '  -0.30%  ' },
    @{ Cell = "D18"; Value = '42.670.11' },
    @{ Cell = "E18"; Value = '  -0.80%  ' },
    @{ Cell = "D19"; Value = '11.98' },
    @{ Cell = "E19"; Value = '  -5.01%  ' },
    @{ Cell = "D20"; Value = '6.14' },
    @{ Cell = "E20"; Value = '  +1.30%  ' },
    @{ Cell = "D21"; Value = '0.0₃0883' },
    @{ Cell = "E21"; Value = '  -1.87%  ' },
    @{ Cell = "D22"; Value = '67.73' },
    @{ Cell = "E22"; Value = '  -0.22%  ' },
    @{ Cell = "D23"; Value = '234.94' },
    @{ Cell = "E23"; Value = '  -0.76%  ' },
    @{ Cell = "D24"; Value = '2.21' },
    @{ Cell = "E24"; Value = '  -0.90%  ' },
    @{ Cell = "E25"; Value = '  +0.00%  ' },
    @{ Cell = "D26"; Value = '2.40' },
    @{ Cell = "E26"; Value = '  -2.11%  ' },
    @{ Cell = "D27"; Value = '24.35' },
    @{ Cell = "E27"; Value = '  -2.29%  ' },
    @{ Cell = "D28"; Value = '2.22' },
    @{ Cell = "E28"; Value = '  +7.44%  ' },
    @{ Cell = "D29"; Value = '9.07' },
    @{ Cell = "E29"; Value = '  -0.80%  ' },
    @{ Cell = "D30"; Value = '31.07' },
    @{ Cell = "E30"; Value = '  -6.67%  ' },
    @{ Cell = "E31"; Value = '  +0.11%  ' },
    @{ Cell = "D32"; Value = '0.0751' },
    @{ Cell = "E32"; Value = '  +7.55%  ' },
    @{ Cell = "D33"; Value = '4.95' },
    @{ Cell = "E33"; Value = '  -1.11%  ' },
    @{ Cell = "B34"; Value = 'Monero' },
    @{ Cell = "C34"; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = "D34"; Value = '132.72' },
    @{ Cell = "E34"; Value = '  -19.74%  ' },
    @{ Cell = "B35"; Value = 'Celestia' },
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia' },
    @{ Cell = "D35"; Value = '17.19' },
    @{ Cell = "E35"; Value = '  -4.58%  ' },
    @{ Cell = "D36"; Value = '2.32' },
    @{ Cell = "E36"; Value = '  -0.98%  ' },
    @{ Cell = "D37"; Value = '1.80' },
    @{ Cell = "E37"; Value = '  +1.96%  ' },
    @{ Cell = "D38"; Value = '4.31' },
    @{ Cell = "E38"; Value = '  -5.09%  ' },
    @{ Cell = "E39"; Value = '  -1.17%  ' },
    @{ Cell = "D40"; Value = '22.16' },
    @{ Cell = "E40"; Value = '  +23.09%  ' },
    @{ Cell = "D41"; Value = '2.72' },
    @{ Cell = "E41"; Value = '  -2.67%  ' },
    @{ Cell = "E42"; Value = '  -1.80%  ' },
    @{ Cell = "D43"; Value = '1.920.65' },
    @{ Cell = "E43"; Value = '  -3.50%  ' },
    @{ Cell = "E44"; Value = '  -0.05%  ' },
    @{ Cell = "D45"; Value = '10.08' },
    @{ Cell = "E45"; Value = '  -5.84%  ' },
    @{ Cell = "E46"; Value = '  +0.05%  ' },
    @{ Cell = "D47"; Value = '2.70' },
    @{ Cell = "E47"; Value = '  -2.71%  ' },
    @{ Cell = "E48"; Value = '  -0.34%  ' },
    @{ Cell = "D49"; Value = '2.553.11' },
    @{ Cell = "E49"; Value = '  +0.03%  ' },
    @{ Cell = "D50"; Value = '52.35' },
    @{ Cell = "E50"; Value = '  -2.71%  ' },
    @{ Cell = "D51"; Value = '71.54' },
    @{ Cell = "E51"; Value = '  -0.76%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $val = $u.Value

    if ($val -match '^\s*[+-]?\d+(\.\d+)?\s*$') {
        $cell.Value = "'" + $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}
